$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3587
$ws.Range("E2").Value = 78
$ws.Range("F2").Value = 78
$ws.Range("G2").Value = 63
$ws.Range("H2").Value = 184
$ws.Range("I2").Value = 184
$ws.Range("K2").Value = 948
$ws.Range("L2").Value = 492
$ws.Range("M2").Value = 456
$ws.Range("N2").Value = 456
$ws.Range("P2").Value = 60
$ws.Range("Q2").Value = 169
$ws.Range("R2").Value = -58
$ws.Range("S2").Value = -120
$ws.Range("T2").Value = 29
$ws.Range("U2").Value = 140
$ws.Range("V2").Value = 120
$ws.Range("W2").Value = 2.18
$ws.Range("X2").Value = 5.12
$ws.Range("Y2").Value = 49.97
$ws.Range("Z2").Value = 20.36
$ws.Range("AA2").Value = 108.01
$ws.Range("AB2").Value = 659.66
$ws.Range("AC2").Value = 1529
$ws.Range("AD2").Value = 2.52
$ws.Range("AE2").Value = 3800
$ws.Range("AF2").Value = 1.02
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1.3
$ws.Range("AI2").Value = 3.27
$ws.Range("AJ2").Value = 12000000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 4367
$ws.Range("E3").Value = 110
$ws.Range("F3").Value = 110
$ws.Range("G3").Value = 108
$ws.Range("H3").Value = 79
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1317
$ws.Range("L3").Value = 783
$ws.Range("M3").Value = 534
$ws.Range("N3").Value = 532
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 60
$ws.Range("Q3").Value = 51
$ws.Range("R3").Value = -223
$ws.Range("S3").Value = 179
$ws.Range("T3").Value = 222
$ws.Range("U3").Value = -171
$ws.Range("V3").Value = 303
$ws.Range("W3").Value = 2.51
$ws.Range("X3").Value = 1.82
$ws.Range("Y3").Value = 16.11
$ws.Range("Z3").Value = 7.02
$ws.Range("AA3").Value = 146.71
$ws.Range("AB3").Value = 785.42
$ws.Range("AC3").Value = 663
$ws.Range("AD3").Value = 18.1
$ws.Range("AE3").Value = 4432
$ws.Range("AF3").Value = 2.71
$ws.Range("AG3").Value = 70
$ws.Range("AH3").Value = 0.58
$ws.Range("AI3").Value = 10.56
$ws.Range("AJ3").Value = 12000000

# Row 4
$ws.Range("D4").Value = 4679
$ws.Range("E4").Value = 89
$ws.Range("F4").Value = 89
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 54
$ws.Range("I4").Value = 55
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 1468
$ws.Range("L4").Value = 895
$ws.Range("M4").Value = 573
$ws.Range("N4").Value = 572
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 60
$ws.Range("Q4").Value = 38
$ws.Range("R4").Value = -181
$ws.Range("S4").Value = 124
$ws.Range("T4").Value = 148
$ws.Range("U4").Value = -110
$ws.Range("V4").Value = 436
$ws.Range("W4").Value = 1.91
$ws.Range("X4").Value = 1.16
$ws.Range("Y4").Value = 9.93
$ws.Range("Z4").Value = 3.89
$ws.Range("AA4").Value = 156.34
$ws.Range("AB4").Value = 853.04
$ws.Range("AC4").Value = 457
$ws.Range("AD4").Value = 16.37
$ws.Range("AE4").Value = 4764
$ws.Range("AF4").Value = 1.57
$ws.Range("AG4").Value = 70
$ws.Range("AH4").Value = 0.9399999999999999
$ws.Range("AI4").Value = 15.32
$ws.Range("AJ4").Value = 12000000

# Row 5
$ws.Range("D5").Value = 5279
$ws.Range("E5").Value = 101
$ws.Range("F5").Value = 101
$ws.Range("G5").Value = 89
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 67
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1882
$ws.Range("L5").Value = 1252
$ws.Range("M5").Value = 630
$ws.Range("N5").Value = 629
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 60
$ws.Range("Q5").Value = 17
$ws.Range("R5").Value = -300
$ws.Range("S5").Value = 289
$ws.Range("T5").Value = 277
$ws.Range("U5").Value = -261
$ws.Range("V5").Value = 733
$ws.Range("W5").Value = 1.92
$ws.Range("X5").Value = 1.27
$ws.Range("Y5").Value = 11.15
$ws.Range("Z5").Value = 3.99
$ws.Range("AA5").Value = 198.59
$ws.Range("AB5").Value = 953.73
$ws.Range("AC5").Value = 558
$ws.Range("AD5").Value = 7.96
$ws.Range("AE5").Value = 5242
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 80
$ws.Range("AH5").Value = 1.8
$ws.Range("AI5").Value = 14.34
$ws.Range("AJ5").Value = 12000000

# Row 6
$ws.Range("D6").Value = 5658
$ws.Range("E6").Value = 89
$ws.Range("F6").Value = 89
$ws.Range("G6").Value = 73
$ws.Range("H6").Value = 53
$ws.Range("I6").Value = 53
$ws.Range("K6").Value = 1996
$ws.Range("L6").Value = 1326
$ws.Range("M6").Value = 670
$ws.Range("N6").Value = 669
$ws.Range("P6").Value = 60
$ws.Range("Q6").Value = 65
$ws.Range("R6").Value = -118
$ws.Range("S6").Value = 47
$ws.Range("T6").Value = 104
$ws.Range("U6").Value = -39
$ws.Range("V6").Value = 790
$ws.Range("W6").Value = 1.58
$ws.Range("X6").Value = 0.93
$ws.Range("Y6").Value = 8.16
$ws.Range("Z6").Value = 2.71
$ws.Range("AA6").Value = 198.01
$ws.Range("AB6").Value = 1019.08
$ws.Range("AC6").Value = 441
$ws.Range("AD6").Value = 9.85
$ws.Range("AE6").Value = 5573
$ws.Range("AF6").Value = 0.78
$ws.Range("AG6").Value = 60
$ws.Range("AH6").Value = 1.38
$ws.Range("AI6").Value = 13.6
$ws.Range("AJ6").Value = 12000000

# Rows 7-9: clear all data columns (D:AJ), keep A/B/C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
